$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.276.23"
$ws.Range("E2").Value = "  +0.72%  "

$ws.Range("D3").Value = "3.343.62"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.34"
$ws.Range("E5").Value = "  +0.30%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.41"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  +1.22%  "

$ws.Range("D8").Value = "3.338.04"
$ws.Range("E8").Value = "  +0.37%  "

$ws.Range("E9").Value = "  -0.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +6.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.634"
$ws.Range("E11").Value = "  +2.38%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.98"
$ws.Range("E12").Value = "  +0.20%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000277"
$ws.Range("E13").Value = "  +2.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.07"
$ws.Range("E14").Value = "  +1.00%  "

$ws.Range("D15").Value = "3.883.12"
$ws.Range("E15").Value = "  +0.25%  "

$ws.Range("E16").Value = "  +2.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.18"
$ws.Range("E17").Value = "  -0.63%  "

$ws.Range("D18").Value = "3.349.57"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").Value = "64.958.59"
$ws.Range("E19").Value = "  +1.84%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.75"
$ws.Range("E20").Value = "  -0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.988"
$ws.Range("E21").Value = "  +1.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "455.02"
$ws.Range("E22").Value = "  +5.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.96"
$ws.Range("E23").Value = "  +7.91%  "

$ws.Range("E24").Value = "  -0.18%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.05"
$ws.Range("E25").Value = "  +7.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "87.04"
$ws.Range("E26").Value = "  +3.76%  "

$ws.Range("E27").Value = "  +2.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.77"
$ws.Range("E28").Value = "  +1.67%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.68"
$ws.Range("E29").Value = "  -0.25%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.98"
$ws.Range("E30").Value = "  +5.13%  "

$ws.Range("E31").Value = "  +1.79%  "

$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "568.01"
$ws.Range("E33").Value = "  -3.77%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.00"
$ws.Range("E34").Value = "  +4.30%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.58"
$ws.Range("E37").Value = "  +3.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.140"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.41"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "0.0₃0739"
$ws.Range("E40").Value = "  -1.11%  "

$ws.Range("E41").Value = "  +1.25%  "

$ws.Range("D42").Value = "3.065.59"
$ws.Range("E42").Value = "  -1.35%  "

$ws.Range("E43").Value = "  -1.26%  "

$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("E45").Value = "  +1.62%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.133"
$ws.Range("E46").Value = "  +3.63%  "

$ws.Range("E47").Value = "  +0.63%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.26"
$ws.Range("E49").Value = "  +4.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.53"
$ws.Range("E50").Value = "  -2.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.16"
$ws.Range("E51").Value = "  +0.08%  "
